$wb = $excel.ActiveWorkbook

# Update "想去人数" (number interested) counts on both the "展览" and
# "全部类型" sheets, which mirror the same event rows.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F6").Value = 457
    $ws.Range("F9").Value = 580
}
